# Auto-generated Excel COM-interop script to apply the diff changes
# to Sheets/Exodus_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1107.1904
$ws.Range("I19").Value = 915.6667
$ws.Range("J19").Value = 1183.8
$ws.Range("K19").Value = 915.6667
$ws.Range("L19").Value = 1183.8
$ws.Range("M19").Value = -740.6667
$ws.Range("N19").Value = -1533.8
$ws.Range("H34").Value = 4450
$ws.Range("I34").Value = 4450
$ws.Range("K34").Value = 4450
$ws.Range("M34").Value = -4247
$ws.Range("H36").Value = 4450
$ws.Range("I36").Value = 4450
$ws.Range("K36").Value = 4450
$ws.Range("M36").Value = -3735
$ws.Range("H40").Value = 9959.700000000001
$ws.Range("I40").Value = 3374
$ws.Range("K40").Value = 3374
$ws.Range("M40").Value = -3199
$ws.Range("H55").Value = 349.1111
$ws.Range("I55").Value = 448
$ws.Range("J55").Value = 225.5
$ws.Range("K55").Value = 448
$ws.Range("L55").Value = 225.5
$ws.Range("M55").Value = -234
$ws.Range("N55").Value = -653.5
$ws.Range("H87").Value = 113749
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 113749
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 113749
$ws.Range("M87").Value = $null
$ws.Range("N87").Value = -116245
$ws.Range("H90").Value = 113749
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 113749
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 341247
$ws.Range("M90").Value = $null
$ws.Range("N90").Value = -353727
$ws.Range("H101").Value = 199037.25
$ws.Range("I101").Value = 3168
$ws.Range("J101").Value = 264327
$ws.Range("K101").Value = 9504
$ws.Range("L101").Value = 792981
$ws.Range("M101").Value = -7882
$ws.Range("N101").Value = -796225
$ws.Range("H103").Value = 1174.5
$ws.Range("I103").Value = 999
$ws.Range("K103").Value = 2997
$ws.Range("M103").Value = -2411
$ws.Range("H107").Value = 1096.3889
$ws.Range("I107").Value = 826.6923
$ws.Range("K107").Value = 826.6923
$ws.Range("M107").Value = 1093.3077
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254
$ws.Range("H121").Value = 1655.6666
$ws.Range("J121").Value = 1655.6666
$ws.Range("L121").Value = 4966.9998
$ws.Range("N121").Value = -8460.9998
$ws.Range("H125").Value = 2998.5
$ws.Range("I125").Value = 2998.5
$ws.Range("K125").Value = 26986.5
$ws.Range("M125").Value = -24526.5
$ws.Range("H133").Value = 98568.42999999999
$ws.Range("J133").Value = 98568.42999999999
$ws.Range("L133").Value = 98568.42999999999
$ws.Range("N133").Value = -108688.43
$ws.Range("H134").Value = 82648.44500000001
$ws.Range("J134").Value = 82648.44500000001
$ws.Range("L134").Value = 82648.44500000001
$ws.Range("N134").Value = -92788.44500000001
$ws.Range("H136").Value = 78605
$ws.Range("J136").Value = 78605
$ws.Range("L136").Value = 78605
$ws.Range("N136").Value = -88805
$ws.Range("H137").Value = 404849.72
$ws.Range("I137").Value = 1323.56
$ws.Range("J137").Value = 1321954.6
$ws.Range("K137").Value = 3970.68
$ws.Range("L137").Value = 3965863.8
$ws.Range("M137").Value = -1420.68
$ws.Range("N137").Value = -3970963.8
$ws.Range("H138").Value = 52686584
$ws.Range("I138").Value = 101483.2
$ws.Range("J138").Value = 111114470
$ws.Range("K138").Value = 304449.6
$ws.Range("L138").Value = 333343410
$ws.Range("M138").Value = -299309.6
$ws.Range("N138").Value = -333353690
$ws.Range("H139").Value = 69518
$ws.Range("J139").Value = 69518
$ws.Range("L139").Value = 69518
$ws.Range("N139").Value = -79798
$ws.Range("H140").Value = 91557.10000000001
$ws.Range("J140").Value = 91557.10000000001
$ws.Range("L140").Value = 91557.10000000001
$ws.Range("N140").Value = -101917.1
$ws.Range("H141").Value = 6448.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 6448.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 19345.5
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = -29705.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7031.271
$ws.Range("I32").Value = 3623.5625
$ws.Range("K32").Value = 3623.5625
$ws.Range("M32").Value = -3336.5625
$ws.Range("H45").Value = 11462869
$ws.Range("I45").Value = 3349.8333
$ws.Range("K45").Value = 3349.8333
$ws.Range("M45").Value = -2972.8333
$ws.Range("H61").Value = 253502.5
$ws.Range("I61").Value = 5005.5
$ws.Range("J61").Value = 501999.5
$ws.Range("K61").Value = 5005.5
$ws.Range("L61").Value = 501999.5
$ws.Range("M61").Value = -4793.5
$ws.Range("N61").Value = -502423.5
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248
$ws.Range("H63").Value = 4329.5
$ws.Range("I63").Value = 3994.25
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 3994.25
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -3308.25
$ws.Range("N63").Value = -6372
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240
$ws.Range("H66").Value = 4329.5
$ws.Range("I66").Value = 3994.25
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 19971.25
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -16539.25
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 4702.1577
$ws.Range("I74").Value = 2381.5
$ws.Range("J74").Value = 11200
$ws.Range("K74").Value = 2381.5
$ws.Range("L74").Value = 11200
$ws.Range("M74").Value = -1507.5
$ws.Range("N74").Value = -12948
$ws.Range("H77").Value = 4702.1577
$ws.Range("I77").Value = 2381.5
$ws.Range("J77").Value = 11200
$ws.Range("K77").Value = 11907.5
$ws.Range("L77").Value = 56000
$ws.Range("M77").Value = -7539.5
$ws.Range("N77").Value = -64736
$ws.Range("H88").Value = 88321
$ws.Range("I88").Value = 825
$ws.Range("J88").Value = 117486.336
$ws.Range("K88").Value = 825
$ws.Range("L88").Value = 117486.336
$ws.Range("M88").Value = -419
$ws.Range("N88").Value = -118298.336
$ws.Range("H91").Value = 88321
$ws.Range("I91").Value = 825
$ws.Range("J91").Value = 117486.336
$ws.Range("K91").Value = 825
$ws.Range("L91").Value = 117486.336
$ws.Range("M91").Value = 579
$ws.Range("N91").Value = -120294.336
$ws.Range("H97").Value = 1321.25
$ws.Range("I97").Value = 795.7143
$ws.Range("K97").Value = 795.7143
$ws.Range("M97").Value = -299.7143
$ws.Range("H104").Value = 79978.664
$ws.Range("J104").Value = 79978.664
$ws.Range("L104").Value = 79978.664
$ws.Range("N104").Value = -86966.664
$ws.Range("H110").Value = 1349.6666
$ws.Range("I110").Value = 989.375
$ws.Range("K110").Value = 989.375
$ws.Range("M110").Value = 1055.625
$ws.Range("H132").Value = 4818.579
$ws.Range("I132").Value = 3962.923
$ws.Range("J132").Value = 6672.5
$ws.Range("K132").Value = 11888.769
$ws.Range("L132").Value = 20017.5
$ws.Range("M132").Value = -9358.769
$ws.Range("N132").Value = -25077.5
$ws.Range("H136").Value = 253502.5
$ws.Range("I136").Value = 5005.5
$ws.Range("J136").Value = 501999.5
$ws.Range("K136").Value = 15016.5
$ws.Range("L136").Value = 1505998.5
$ws.Range("M136").Value = -12466.5
$ws.Range("N136").Value = -1511098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 40500
$ws.Range("I26").Value = 40500
$ws.Range("K26").Value = 40500
$ws.Range("M26").Value = -40208
$ws.Range("H86").Value = 6394.737
$ws.Range("I86").Value = 3717.2856
$ws.Range("J86").Value = 7956.5835
$ws.Range("K86").Value = 3717.2856
$ws.Range("L86").Value = 7956.5835
$ws.Range("M86").Value = -2594.2856
$ws.Range("N86").Value = -10202.5835
$ws.Range("H87").Value = 175000
$ws.Range("J87").Value = 175000
$ws.Range("L87").Value = 175000
$ws.Range("N87").Value = -177496
$ws.Range("H89").Value = 6394.737
$ws.Range("I89").Value = 3717.2856
$ws.Range("J89").Value = 7956.5835
$ws.Range("K89").Value = 18586.428
$ws.Range("L89").Value = 39782.9175
$ws.Range("M89").Value = -12970.428
$ws.Range("N89").Value = -51014.9175
$ws.Range("H90").Value = 175000
$ws.Range("J90").Value = 175000
$ws.Range("L90").Value = 525000
$ws.Range("N90").Value = -537480
$ws.Range("H94").Value = 1727.4667
$ws.Range("I94").Value = 1624.1538
$ws.Range("K94").Value = 1624.1538
$ws.Range("M94").Value = -1173.1538
$ws.Range("H99").Value = 44693.75
$ws.Range("I99").Value = 61444.35
$ws.Range("J99").Value = 4013.7144
$ws.Range("K99").Value = 61444.35
$ws.Range("L99").Value = 4013.7144
$ws.Range("M99").Value = -59946.35
$ws.Range("N99").Value = -7009.7144
$ws.Range("H107").Value = 2725.4614
$ws.Range("I107").Value = 2744.25
$ws.Range("K107").Value = 2744.25
$ws.Range("M107").Value = -824.25
$ws.Range("H135").Value = 97665.71000000001
$ws.Range("J135").Value = 97665.71000000001
$ws.Range("L135").Value = 97665.71000000001
$ws.Range("N135").Value = -107805.71
$ws.Range("H138").Value = 82370.375
$ws.Range("J138").Value = 82370.375
$ws.Range("L138").Value = 82370.375
$ws.Range("N138").Value = -92650.375
$ws.Range("H140").Value = 75711.28999999999
$ws.Range("J140").Value = 75711.28999999999
$ws.Range("L140").Value = 75711.28999999999
$ws.Range("N140").Value = -86071.28999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1876
$ws.Range("I22").Value = 1851.2
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1851.2
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1501.2
$ws.Range("N22").Value = -2700
$ws.Range("H31").Value = 3098.8235
$ws.Range("I31").Value = 2029.875
$ws.Range("J31").Value = 4049
$ws.Range("K31").Value = 2029.875
$ws.Range("L31").Value = 4049
$ws.Range("M31").Value = -1734.875
$ws.Range("N31").Value = -4639
$ws.Range("H34").Value = 3098.8235
$ws.Range("I34").Value = 2029.875
$ws.Range("J34").Value = 4049
$ws.Range("K34").Value = 2029.875
$ws.Range("L34").Value = 4049
$ws.Range("M34").Value = -1827.875
$ws.Range("N34").Value = -4453
$ws.Range("H58").Value = 2805.8333
$ws.Range("J58").Value = 3406.25
$ws.Range("L58").Value = 3406.25
$ws.Range("N58").Value = -3812.25
$ws.Range("H86").Value = 3263861.8
$ws.Range("J86").Value = 21698.8
$ws.Range("L86").Value = 21698.8
$ws.Range("N86").Value = -23944.8
$ws.Range("H89").Value = 3263861.8
$ws.Range("J89").Value = 21698.8
$ws.Range("L89").Value = 108494
$ws.Range("N89").Value = -119726
$ws.Range("H92").Value = 54666.668
$ws.Range("J92").Value = 54666.668
$ws.Range("L92").Value = 54666.668
$ws.Range("N92").Value = -59658.668
$ws.Range("H93").Value = 10171
$ws.Range("I93").Value = 10171
$ws.Range("K93").Value = 10171
$ws.Range("M93").Value = -8299
$ws.Range("H105").Value = 2991.4736
$ws.Range("I105").Value = 1194.75
$ws.Range("K105").Value = 1194.75
$ws.Range("M105").Value = 552.25
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null
$ws.Range("H122").Value = 3874.65
$ws.Range("I122").Value = 3715.5715
$ws.Range("J122").Value = 3960.3076
$ws.Range("K122").Value = 11146.7145
$ws.Range("L122").Value = 11880.9228
$ws.Range("M122").Value = -8696.7145
$ws.Range("N122").Value = -16780.9228
$ws.Range("H132").Value = 2788906.8
$ws.Range("I132").Value = 2847127
$ws.Range("J132").Value = 2602602.2
$ws.Range("K132").Value = 8541381
$ws.Range("L132").Value = 7807806.600000001
$ws.Range("M132").Value = -8538851
$ws.Range("N132").Value = -7812866.600000001
$ws.Range("H134").Value = 3249759.5
$ws.Range("I134").Value = 4764914
$ws.Range("K134").Value = 14294742
$ws.Range("M134").Value = -14292207
$ws.Range("H136").Value = 2805.8333
$ws.Range("J136").Value = 3406.25
$ws.Range("L136").Value = 10218.75
$ws.Range("N136").Value = -15318.75
$ws.Range("H138").Value = 92455.625
$ws.Range("J138").Value = 92455.625
$ws.Range("L138").Value = 92455.625
$ws.Range("N138").Value = -102735.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 89.8
$ws.Range("I2").Value = 89.8
$ws.Range("K2").Value = 538.8
$ws.Range("M2").Value = -425.8
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H131").Value = 1442.0834
$ws.Range("I131").Value = 1149.1666
$ws.Range("J131").Value = 1735
$ws.Range("K131").Value = 3447.4998
$ws.Range("L131").Value = 5205
$ws.Range("M131").Value = 1592.5002
$ws.Range("N131").Value = -15285
$ws.Range("H138").Value = 2218.5715
$ws.Range("I138").Value = 1817.1875
$ws.Range("K138").Value = 5451.5625
$ws.Range("M138").Value = -311.5625
$ws.Range("H139").Value = 1457.6666
$ws.Range("I139").Value = 1457.6666
$ws.Range("K139").Value = 4372.9998
$ws.Range("M139").Value = 767.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2016.963
$ws.Range("I97").Value = 2231.75
$ws.Range("J97").Value = 1704.5454
$ws.Range("K97").Value = 2231.75
$ws.Range("L97").Value = 1704.5454
$ws.Range("M97").Value = -1735.75
$ws.Range("N97").Value = -2696.5454
$ws.Range("H102").Value = 1361.75
$ws.Range("I102").Value = 1396.9429
$ws.Range("K102").Value = 1396.9429
$ws.Range("M102").Value = 225.0571
$ws.Range("H109").Value = 86995.664
$ws.Range("J109").Value = 86995.664
$ws.Range("L109").Value = 86995.664
$ws.Range("N109").Value = -89075.664
$ws.Range("H122").Value = 5007003
$ws.Range("I122").Value = 5007003
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15021009
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15018559
$ws.Range("N122").Value = $null
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H135").Value = 99916.44
$ws.Range("J135").Value = 99916.44
$ws.Range("L135").Value = 99916.44
$ws.Range("N135").Value = -110056.44
$ws.Range("H140").Value = 90251.75
$ws.Range("J140").Value = 90251.75
$ws.Range("L140").Value = 90251.75
$ws.Range("N140").Value = -100611.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1335.1428
$ws.Range("J16").Value = 2271.25
$ws.Range("L16").Value = 2271.25
$ws.Range("N16").Value = -2611.25
$ws.Range("H46").Value = 4200.4
$ws.Range("I46").Value = 2750.5
$ws.Range("J46").Value = 5167
$ws.Range("K46").Value = 2750.5
$ws.Range("L46").Value = 5167
$ws.Range("M46").Value = -2562.5
$ws.Range("N46").Value = -5543
$ws.Range("H55").Value = 3750
$ws.Range("I55").Value = 2562.625
$ws.Range("J55").Value = 8499.5
$ws.Range("K55").Value = 2562.625
$ws.Range("L55").Value = 8499.5
$ws.Range("M55").Value = -2389.625
$ws.Range("N55").Value = -8845.5
$ws.Range("H93").Value = 3510.6
$ws.Range("I93").Value = 3137.75
$ws.Range("J93").Value = 5002
$ws.Range("K93").Value = 3137.75
$ws.Range("L93").Value = 5002
$ws.Range("M93").Value = -1889.75
$ws.Range("N93").Value = -7498
$ws.Range("H100").Value = 4411.1055
$ws.Range("I100").Value = 3363.3125
$ws.Range("J100").Value = 9999.333000000001
$ws.Range("K100").Value = 3363.3125
$ws.Range("L100").Value = 9999.333000000001
$ws.Range("M100").Value = -2822.3125
$ws.Range("N100").Value = -11081.333
$ws.Range("J122").Value = 6426.3335
$ws.Range("L122").Value = 19279.0005
$ws.Range("N122").Value = -24179.0005
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null
$ws.Range("H132").Value = 2778.2222
$ws.Range("I132").Value = 1917.8334
$ws.Range("K132").Value = 5753.5002
$ws.Range("M132").Value = -3223.5002
$ws.Range("H136").Value = 4725.5835
$ws.Range("I136").Value = 4362.4
$ws.Range("K136").Value = 13087.2
$ws.Range("M136").Value = -10537.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 27000
$ws.Range("J63").Value = 27000
$ws.Range("L63").Value = 27000
$ws.Range("N63").Value = -28248
$ws.Range("H66").Value = 27000
$ws.Range("J66").Value = 27000
$ws.Range("L66").Value = 81000
$ws.Range("N66").Value = -87240
$ws.Range("H74").Value = 32899.8
$ws.Range("J74").Value = 32899.8
$ws.Range("L74").Value = 32899.8
$ws.Range("N74").Value = -34771.8
$ws.Range("H77").Value = 32899.8
$ws.Range("J77").Value = 32899.8
$ws.Range("L77").Value = 98699.40000000001
$ws.Range("N77").Value = -108059.4
$ws.Range("H107").Value = 1249.2858
$ws.Range("I107").Value = 1502.6666
$ws.Range("K107").Value = 4507.9998
$ws.Range("M107").Value = -2587.9998
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800
$ws.Range("H132").Value = 3075.9048
$ws.Range("I132").Value = 2808
$ws.Range("J132").Value = 3511.25
$ws.Range("K132").Value = 8424
$ws.Range("L132").Value = 10533.75
$ws.Range("M132").Value = -5894
$ws.Range("N132").Value = -15593.75
$ws.Range("H136").Value = 1938.1875
$ws.Range("I136").Value = 1323.7273
$ws.Range("K136").Value = 3971.1819
$ws.Range("M136").Value = -1421.1819
